# Applies the "best current data" update described in the commit message:
#  - refreshed measurement inputs for the CONTROL / STATIC / DYNAMIC sample tables
#  - row 10 (STATIC / Sample 2) measurements removed -> formulas fall back to errors/0
#  - row 9 (STATIC / Sample 1) gains an Width->E9 measurement it was previously missing
#  - UI state: active selection moves to I7, window is minimized/shifted

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (CONTROL, Sample 1): updated measurements ---
$ws.Range("B3").Value = 40
$ws.Range("C3").Value = 17.309999999999999
$ws.Range("D3").Value = 1.58
$ws.Range("E3").Value = 1.1200000000000001

# --- Row 5: G5 was a stale literal left over from before the fitting code was
#     finished; replace it with the live formula (matches the G4:G5 shared pattern) ---
$ws.Range("G5").Formula = "=(PI()/4)*D5*E5"

# --- Row 9 (STATIC, Sample 1): updated measurements + newly recorded width (E9) ---
$ws.Range("B9").Value = 103.23
$ws.Range("C9").Value = 15.57
$ws.Range("D9").Value = 3.6
$ws.Range("E9").Value = 1.05

# --- Row 10 (STATIC, Sample 2): measurements removed/unavailable this round ---
$ws.Range("B10:E10").ClearContents()

# --- Row 15 (DYNAMIC, Sample 1): updated measurements ---
$ws.Range("B15").Value = 65
$ws.Range("C15").Value = 20.74
$ws.Range("D15").Value = 3.91
$ws.Range("E15").Value = 0.75

# --- Active selection moves from B15 to I7 ---
$ws.Range("I7").Select()

# --- Window state: minimized, shifted horizontally ---
$win = $excel.ActiveWindow
$win.WindowState = -4140
$win.Left = 3960

$wb.Saved = $false
